$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previously entered shoot-growth data while keeping cell formatting.
# Rows 2-25: columns D through U held measurements/formulas that are being revised out.
$ws.Range("D2:U25").ClearContents()

# Rows 26-30: column C (bag/collector label) through U also cleared.
$ws.Range("C26:U30").ClearContents()

# Update the active selection to D1, matching the saved view state.
$ws.Range("D1").Select()
